# Applies the "rendering fixes + start working on movement" edit described
# by the OOXML diff to the currently-open workbook (resources/map.xlsx).
#
# Summary of the change:
#  1. The grid (sheet "Blad1") is extended with 10 extra columns (X:AG) on
#     every data row (2-37), all filled with "GROUND".
#  2. The "shelf" object that used to sit at rows 4-5 / columns C-J is
#     moved to rows 7-8 / columns E-L (2 columns to the right, 3 rows down).
#     The cells it vacates become plain "GROUND" again.
#  3. While being moved, the shelf's shared strings are renamed:
#       SHELF_N_CC:aa -> SHELF_or:N_cc:aa
#       SHELF_N_CC:ab -> SHELF_or:N_cc:ab
#       SHELF_S_CC:ac -> SHELF_or:S_cc:ac
#       SHELF_S_CC:ad -> SHELF_or:S_cc:ad
#  4. Cosmetic view changes: the sheet is zoomed to 55% and the active
#     selection moves from L18 to J16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill the newly-used columns X:AG (rows 2-37) with GROUND --------
$ws.Range("X2:AG37").Value = "GROUND"

# --- 2 & 3. Move the shelf block from C4:J5 to E7:L8 ---------------------
# Vacate the old location, turning it back into plain floor.
$ws.Range("C4:J5").Value = "GROUND"

# Place the (renamed) shelf pieces in their new location.
$ws.Range("E7:H7").Value = "SHELF_or:N_cc:aa"
$ws.Range("I7:L7").Value = "SHELF_or:N_cc:ab"
$ws.Range("E8:H8").Value = "SHELF_or:S_cc:ac"
$ws.Range("I8:L8").Value = "SHELF_or:S_cc:ad"

# --- 4. View/selection cosmetics -----------------------------------------
$excel.ActiveWindow.Zoom = 55
$ws.Range("J16").Select() | Out-Null
